$d = $word.ActiveDocument

function Get-ParagraphByText($doc, $needle) {
    foreach ($p in $doc.Paragraphs) {
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $needle) {
            return $p
        }
    }
    return $null
}

# --- 1. Add the new "autotiling" bullet right after the "Towns have
#        affinity..." bullet (same list, same numId/style/font). ---
$anchorPara = Get-ParagraphByText $d "Towns have affinity towards you based on how many tasks you do for residents?"
if ($anchorPara -ne $null) {
    $anchorPara.Range.InsertParagraphAfter() | Out-Null
    $newPara = $d.Paragraphs($anchorPara.Index + 1)
    $newPara.Range.Text = "Towns should have chances to be “different” fundamentally, like some are mining towns, some are merchant towns, some are filled with inns, others with lots of shops, markets, some are fishing towns, others are crafting towns, magic towns, knight strongholds, all kinds of different seeds for the towns so it’s interesting finding new ones."
}

# --- 2. Clear the stale <w:lastRenderedPageBreak/> cached in front of
#        the "Cave dungeon" run. ---
$cavePara = Get-ParagraphByText $d "Cave dungeon"
if ($cavePara -ne $null) {
    $cavePara.Range.Text = "Cave dungeon"
}
